# Updates cryptos list: refreshed prices / 1h volume %, plus a few coins
# that swapped rank position (Toncoin/PancakeSwap, EthereumClassic/Kaspa,
# FTXToken/Maker/Aave). Price values that look like plain decimals are
# written with a leading "'" so Excel keeps them as text (matching the
# sheet's existing text-formatted Price column) instead of coercing them
# to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.324.83"
$ws.Range("E2").Value = "  +2.07%  "

$ws.Range("D3").Value = "2.061.25"
$ws.Range("E3").Value = "  +3.32%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'234.75"
$ws.Range("E5").Value = "  -0.46%  "

$ws.Range("E6").Value = "  +2.27%  "

$ws.Range("D7").Value = "'58.25"
$ws.Range("E7").Value = "  +6.23%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").Value = "'0.382"
$ws.Range("E9").Value = "  +2.79%  "

$ws.Range("D10").Value = "'59.07"
$ws.Range("E10").Value = "  +1.34%  "

$ws.Range("E11").Value = "  +1.82%  "

$ws.Range("E12").Value = "  +2.64%  "

$ws.Range("D13").Value = "2.365.02"
$ws.Range("E13").Value = "  +3.38%  "

$ws.Range("D14").Value = "'14.59"
$ws.Range("E14").Value = "  +2.48%  "

$ws.Range("D15").Value = "'21.19"
$ws.Range("E15").Value = "  +3.75%  "

$ws.Range("D16").Value = "'0.775"
$ws.Range("E16").Value = "  +2.35%  "

$ws.Range("D17").Value = "'5.18"
$ws.Range("E17").Value = "  +2.04%  "

$ws.Range("D18").Value = "2.066.14"
$ws.Range("E18").Value = "  +3.10%  "

$ws.Range("D19").Value = "37.524.66"
$ws.Range("E19").Value = "  +2.78%  "

$ws.Range("D20").Value = "'6.12"
$ws.Range("E20").Value = "  +15.79%  "

$ws.Range("D21").Value = "'70.05"
$ws.Range("E21").Value = "  +3.25%  "

$ws.Range("E22").Value = "  +1.01%  "

$ws.Range("D23").Value = "'227.33"
$ws.Range("E23").Value = "  +2.40%  "

$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "'2.43"
$ws.Range("E25").Value = "  +1.57%  "

$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'2.39"
$ws.Range("E26").Value = "  +1.19%  "

$ws.Range("D27").Value = "'165.39"
$ws.Range("E27").Value = "  +1.92%  "

$ws.Range("E28").Value = "  +11.74%  "

$ws.Range("D29").Value = "'8.86"
$ws.Range("E29").Value = "  +2.31%  "

$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "'0.128"
$ws.Range("E30").Value = "  +0.11%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'19.15"
$ws.Range("E31").Value = "  +1.40%  "

$ws.Range("D32").Value = "'0.119"
$ws.Range("E32").Value = "  +2.14%  "

$ws.Range("D33").Value = "'4.50"
$ws.Range("E33").Value = "  +2.94%  "

$ws.Range("D34").Value = "'0.0622"
$ws.Range("E34").Value = "  +2.61%  "

$ws.Range("E35").Value = "  +9.13%  "

$ws.Range("E36").Value = "  +6.62%  "

$ws.Range("E37").Value = "  +0.66%  "

$ws.Range("E38").Value = "  +0.10%  "

$ws.Range("E39").Value = "  +1.30%  "

$ws.Range("E40").Value = "  +3.67%  "

$ws.Range("D41").Value = "'0.0982"
$ws.Range("E41").Value = "  +3.77%  "

$ws.Range("E42").Value = "  -1.31%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.455.95"
$ws.Range("E43").Value = "  +0.22%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'95.42"
$ws.Range("E44").Value = "  +7.09%  "

$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").Value = "'4.31"
$ws.Range("E45").Value = "  +21.39%  "

$ws.Range("D46").Value = "'0.0211"
$ws.Range("E46").Value = "  +4.24%  "

$ws.Range("D47").Value = "'1.16"
$ws.Range("E47").Value = "  +4.73%  "

$ws.Range("D48").Value = "'15.82"
$ws.Range("E48").Value = "  +3.79%  "

$ws.Range("D49").Value = "'1.03"
$ws.Range("E49").Value = "  +3.39%  "

$ws.Range("D50").Value = "'7.25"
$ws.Range("E50").Value = "  +5.69%  "

$ws.Range("E51").Value = "  +2.04%  "
